$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2116.7307
$ws.Range("J17").Value = 2236.3044
$ws.Range("L17").Value = 6708.9132
$ws.Range("N17").Value = -7044.9132
$ws.Range("H112").Value = 1081.2858
$ws.Range("J112").Value = 1888
$ws.Range("L112").Value = 5664
$ws.Range("N112").Value = -7880
$ws.Range("H127").Value = 3475.125
$ws.Range("I127").Value = 4133.8335
$ws.Range("J127").Value = 1499
$ws.Range("K127").Value = 12401.5005
$ws.Range("L127").Value = 4497
$ws.Range("M127").Value = -7441.500499999998
$ws.Range("N127").Value = -14417
$ws.Range("H132").Value = 35478.1
$ws.Range("I132").Value = 38651.76
$ws.Range("K132").Value = 115955.28
$ws.Range("M132").Value = -113425.28
$ws.Range("H138").Value = 4126.815
$ws.Range("J138").Value = 3702.068
$ws.Range("L138").Value = 11106.204
$ws.Range("N138").Value = -21386.204

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1241.258
$ws.Range("I2").Value = 1136.3922
$ws.Range("J2").Value = 1727.4546
$ws.Range("K2").Value = 1136.3922
$ws.Range("L2").Value = 1727.4546
$ws.Range("M2").Value = -1023.3922
$ws.Range("N2").Value = -1953.4546
$ws.Range("H32").Value = 5159678
$ws.Range("I32").Value = 5380997.5
$ws.Range("J32").Value = 14000
$ws.Range("K32").Value = 5380997.5
$ws.Range("L32").Value = 14000
$ws.Range("M32").Value = -5380710.5
$ws.Range("N32").Value = -14574
$ws.Range("H45").Value = 4646.1577
$ws.Range("I45").Value = 3986.8823
$ws.Range("K45").Value = 3986.8823
$ws.Range("M45").Value = -3609.8823
$ws.Range("H61").Value = 1086514.4
$ws.Range("I61").Value = 1676496.9
$ws.Range("J61").Value = 13819
$ws.Range("K61").Value = 1676496.9
$ws.Range("L61").Value = 13819
$ws.Range("M61").Value = -1676284.9
$ws.Range("N61").Value = -14243
$ws.Range("H74").Value = 3381827.8
$ws.Range("I74").Value = 4169749.5
$ws.Range("J74").Value = 5019.5713
$ws.Range("K74").Value = 4169749.5
$ws.Range("L74").Value = 5019.5713
$ws.Range("M74").Value = -4168875.5
$ws.Range("N74").Value = -6767.5713
$ws.Range("H77").Value = 3381827.8
$ws.Range("I77").Value = 4169749.5
$ws.Range("J77").Value = 5019.5713
$ws.Range("K77").Value = 20848747.5
$ws.Range("L77").Value = 25097.8565
$ws.Range("M77").Value = -20844379.5
$ws.Range("N77").Value = -33833.85649999999
$ws.Range("H88").Value = 3047.5
$ws.Range("J88").Value = 3047.5
$ws.Range("L88").Value = 3047.5
$ws.Range("N88").Value = -3859.5
$ws.Range("H91").Value = 3047.5
$ws.Range("J91").Value = 3047.5
$ws.Range("L91").Value = 3047.5
$ws.Range("N91").Value = -5855.5
$ws.Range("H97").Value = 1122.8334
$ws.Range("I97").Value = 1372.5
$ws.Range("J97").Value = 623.5
$ws.Range("K97").Value = 1372.5
$ws.Range("L97").Value = 623.5
$ws.Range("M97").Value = -876.5
$ws.Range("N97").Value = -1615.5
$ws.Range("H116").Value = 1241.258
$ws.Range("I116").Value = 1136.3922
$ws.Range("J116").Value = 1727.4546
$ws.Range("K116").Value = 1136.3922
$ws.Range("L116").Value = 1727.4546
$ws.Range("M116").Value = 1157.6078
$ws.Range("N116").Value = -6315.4546
$ws.Range("H132").Value = 404516.66
$ws.Range("I132").Value = 569718.5
$ws.Range("K132").Value = 1709155.5
$ws.Range("M132").Value = -1706625.5
$ws.Range("H136").Value = 1086514.4
$ws.Range("I136").Value = 1676496.9
$ws.Range("J136").Value = 13819
$ws.Range("K136").Value = 5029490.699999999
$ws.Range("L136").Value = 41457
$ws.Range("M136").Value = -5026940.699999999
$ws.Range("N136").Value = -46557

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1241.258
$ws.Range("I3").Value = 1136.3922
$ws.Range("J3").Value = 1727.4546
$ws.Range("K3").Value = 1136.3922
$ws.Range("L3").Value = 1727.4546
$ws.Range("M3").Value = -1022.3922
$ws.Range("N3").Value = -1955.4546
$ws.Range("H134").Value = 599538.7
$ws.Range("I134").Value = 885772.4
$ws.Range("J134").Value = 5053.385
$ws.Range("K134").Value = 2657317.2
$ws.Range("L134").Value = 15160.155
$ws.Range("M134").Value = -2654782.2
$ws.Range("N134").Value = -20230.155

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9715.6875
$ws.Range("I31").Value = 3638.5
$ws.Range("K31").Value = 3638.5
$ws.Range("M31").Value = -3343.5
$ws.Range("H34").Value = 9715.6875
$ws.Range("I34").Value = 3638.5
$ws.Range("K34").Value = 3638.5
$ws.Range("M34").Value = -3436.5
$ws.Range("H86").Value = 2421.8125
$ws.Range("I86").Value = 2339.4443
$ws.Range("K86").Value = 2339.4443
$ws.Range("M86").Value = -1216.4443
$ws.Range("H89").Value = 2421.8125
$ws.Range("I89").Value = 2339.4443
$ws.Range("K89").Value = 11697.2215
$ws.Range("M89").Value = -6081.2215
$ws.Range("H105").Value = 1321.963
$ws.Range("I105").Value = 1321.963
$ws.Range("K105").Value = 1321.963
$ws.Range("M105").Value = 425.037
$ws.Range("H134").Value = 2795.6316
$ws.Range("I134").Value = 1918.8235
$ws.Range("K134").Value = 5756.470499999999
$ws.Range("M134").Value = -3221.470499999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4862
$ws.Range("I3").Value = 1297.2307
$ws.Range("K3").Value = 3891.6921
$ws.Range("M3").Value = -3779.6921
$ws.Range("H5").Value = 3562.6316
$ws.Range("I5").Value = 648.5
$ws.Range("J5").Value = 5682
$ws.Range("K5").Value = 1945.5
$ws.Range("L5").Value = 17046
$ws.Range("M5").Value = -1833.5
$ws.Range("N5").Value = -17270
$ws.Range("H68").Value = 6431.6665
$ws.Range("I68").Value = 2456.8
$ws.Range("J68").Value = 7673.8125
$ws.Range("K68").Value = 7370.400000000001
$ws.Range("L68").Value = 23021.4375
$ws.Range("M68").Value = -6559.400000000001
$ws.Range("N68").Value = -24643.4375
$ws.Range("H71").Value = 6431.6665
$ws.Range("I71").Value = 2456.8
$ws.Range("J71").Value = 7673.8125
$ws.Range("K71").Value = 22111.2
$ws.Range("L71").Value = 69064.3125
$ws.Range("M71").Value = -18055.2
$ws.Range("N71").Value = -77176.3125
$ws.Range("H113").Value = 1758.1428
$ws.Range("J113").Value = 1759.16
$ws.Range("L113").Value = 5277.48
$ws.Range("N113").Value = -9617.48
$ws.Range("H131").Value = 15184.277
$ws.Range("I131").Value = 10500
$ws.Range("J131").Value = 15769.8125
$ws.Range("K131").Value = 31500
$ws.Range("L131").Value = 47309.4375
$ws.Range("M131").Value = -26460
$ws.Range("N131").Value = -57389.4375
$ws.Range("H135").Value = 3562.6316
$ws.Range("I135").Value = 648.5
$ws.Range("J135").Value = 5682
$ws.Range("K135").Value = 5836.5
$ws.Range("L135").Value = 51138
$ws.Range("M135").Value = -3301.5
$ws.Range("N135").Value = -56208
$ws.Range("H139").Value = 877.36365
$ws.Range("I139").Value = 765.1
$ws.Range("K139").Value = 2295.3
$ws.Range("M139").Value = 2844.7

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2818.5
$ws.Range("I122").Value = 1607.5483
$ws.Range("J122").Value = 15331.667
$ws.Range("K122").Value = 4822.644899999999
$ws.Range("L122").Value = 45995.001
$ws.Range("M122").Value = -2372.644899999999
$ws.Range("N122").Value = -50895.001
$ws.Range("H132").Value = 233245.84
$ws.Range("I132").Value = 327137.9
$ws.Range("J132").Value = 1645.4
$ws.Range("K132").Value = 981413.7000000001
$ws.Range("L132").Value = 4936.200000000001
$ws.Range("M132").Value = -978883.7000000001
$ws.Range("N132").Value = -9996.200000000001
$ws.Range("H136").Value = 20376.451
$ws.Range("J136").Value = 20376.451
$ws.Range("L136").Value = 61129.353
$ws.Range("N136").Value = -66229.353

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 19626
$ws.Range("I40").Value = 19626
$ws.Range("K40").Value = 19626
$ws.Range("M40").Value = -19490
$ws.Range("H93").Value = 2220.5
$ws.Range("J93").Value = 2269
$ws.Range("L93").Value = 2269
$ws.Range("N93").Value = -4765
$ws.Range("H122").Value = 40880.297
$ws.Range("I122").Value = 3330.4119
$ws.Range("K122").Value = 9991.235700000001
$ws.Range("M122").Value = -7541.235700000001
$ws.Range("H132").Value = 739199
$ws.Range("I132").Value = 963382.0600000001
$ws.Range("J132").Value = 5508.909
$ws.Range("K132").Value = 2890146.18
$ws.Range("L132").Value = 16526.727
$ws.Range("M132").Value = -2887616.18
$ws.Range("N132").Value = -21586.727

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4539.8
$ws.Range("I62").Value = 4174.75
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 4174.75
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -3550.75
$ws.Range("N62").Value = -7248
$ws.Range("H65").Value = 4539.8
$ws.Range("I65").Value = 4174.75
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 20873.75
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -17753.75
$ws.Range("N65").Value = -36240
$ws.Range("H81").Value = 1331.75
$ws.Range("I81").Value = 1053.6666
$ws.Range("J81").Value = 2166
$ws.Range("K81").Value = 2107.3332
$ws.Range("L81").Value = 4332
$ws.Range("M81").Value = -1046.3332
$ws.Range("N81").Value = -6454
$ws.Range("H84").Value = 1331.75
$ws.Range("I84").Value = 1053.6666
$ws.Range("J84").Value = 2166
$ws.Range("K84").Value = 10536.666
$ws.Range("L84").Value = 21660
$ws.Range("M84").Value = -5232.666000000001
$ws.Range("N84").Value = -32268
$ws.Range("H107").Value = 1462.0541
$ws.Range("I107").Value = 1393.5
$ws.Range("K107").Value = 4180.5
$ws.Range("M107").Value = -2260.5
$ws.Range("H136").Value = 8088829.5
$ws.Range("I136").Value = 9051067
$ws.Range("J136").Value = 6038.8
$ws.Range("K136").Value = 27153201
$ws.Range("L136").Value = 18116.4
$ws.Range("M136").Value = -27150651
$ws.Range("N136").Value = -23216.4
